$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "cv"
$ws.Range("C1").Value = "lb"

# Data rows: id values 1..30 in column A, rows 2..31
for ($i = 1; $i -le 30; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $i
}

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

$ws.Range("B2").Select()
